$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regular text cell updates ---
$ws.Range("D1").Value = "Linear SVM"
$ws.Range("E1").Value = "RBF SVM"
$ws.Range("D2").Value = "anxiety_model_20250510_1849"
$ws.Range("E2").Value = "anxiety_model_20250510_1519"
$ws.Range("D3").Value = "I hate everything, I don't want to live anymore."
$ws.Range("E3").Value = "I hate everything, I don't want to live anymore."
$ws.Range("C4").Value = "Anksiyete (84.31%)"
$ws.Range("D4").Value = "Anksiyete (80.24%)"
$ws.Range("E4").Value = "Anksiyete (75.17%)"
$ws.Range("D5").Value = "I feel happy and excited for the day ahead."
$ws.Range("E5").Value = "I feel happy and excited for the day ahead."
$ws.Range("C6").Value = "Normal (14.91%)"
$ws.Range("D6").Value = "Normal (33.03%)"
$ws.Range("E6").Value = "Normal (38.30%)"
$ws.Range("D7").Value = "Life feels like a burden I can’t carry anymore."
$ws.Range("E7").Value = "Life feels like a burden I can’t carry anymore."
$ws.Range("C8").Value = "Anksiyete (69.31%)"
$ws.Range("D8").Value = "Anksiyete (85.16%)"
$ws.Range("E8").Value = "Anksiyete (78.01%)"
$ws.Range("D9").Value = "I enjoyed a lovely walk in the park today."
$ws.Range("E9").Value = "I enjoyed a lovely walk in the park today."
$ws.Range("C10").Value = "Normal (5.23%)"
$ws.Range("D10").Value = "Normal (18.81%)"
$ws.Range("E10").Value = "Normal (26.02%)"
$ws.Range("D11").Value = "I just want to disappear and never come back."
$ws.Range("E11").Value = "I just want to disappear and never come back."
$ws.Range("C12").Value = "Normal (30.06%)"
$ws.Range("D12").Value = "Anksiyete (58.30%)"
$ws.Range("E12").Value = "Anksiyete (55.40%)"
$ws.Range("D13").Value = "I’m looking forward to spending time with my friends."
$ws.Range("E13").Value = "I’m looking forward to spending time with my friends."
$ws.Range("C14").Value = "Normal (24.60%)"
$ws.Range("E14").Value = "Normal (48.54%)"
$ws.Range("D15").Value = "Nothing I do seems to matter; I feel so empty."
$ws.Range("E15").Value = "Nothing I do seems to matter; I feel so empty."
$ws.Range("C16").Value = "Normal (41.52%)"
$ws.Range("D16").Value = "Anksiyete (67.67%)"
$ws.Range("E16").Value = "Anksiyete (59.60%)"
$ws.Range("D17").Value = "I’m grateful for the little things that make me smile."
$ws.Range("E17").Value = "I’m grateful for the little things that make me smile."
$ws.Range("C18").Value = "Normal (8.27%)"
$ws.Range("D18").Value = "Normal (36.99%)"
$ws.Range("E18").Value = "Normal (35.28%)"
$ws.Range("D19").Value = "Even surrounded by people, I feel completely alone."
$ws.Range("E19").Value = "Even surrounded by people, I feel completely alone."
$ws.Range("C20").Value = "Normal (38.21%)"
$ws.Range("D20").Value = "Anksiyete (52.81%)"
$ws.Range("E20").Value = "Anksiyete (54.88%)"
$ws.Range("D21").Value = "Today was a productive and fulfilling day."
$ws.Range("E21").Value = "Today was a productive and fulfilling day."
$ws.Range("C22").Value = "Normal (5.13%)"
$ws.Range("D22").Value = "Normal (17.71%)"
$ws.Range("E22").Value = "Normal (25.22%)"
$ws.Range("A24").Value = "CROSS VALIDATION  ACCURACY:"

# --- D14 needs its style flipped from "Normal"(green) to "Anksiyete"(red);
#     copy format from D16 (already red-styled) then set the text value ---
$ws.Range("D16").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "Anksiyete (55.79%)"

# --- Numeric-looking text cells (must remain stored as text, not numbers) ---
# Assign with a leading apostrophe to force text, then repaste formats from the
# untouched neighbouring "F" cell in the same row to keep the original cell style.
$ws.Range("D24").Value = "'0.9188"
$ws.Range("F24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "'0.9202"
$ws.Range("F24").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("D25").Value = "'0.9251"
$ws.Range("F25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "'0.9254"
$ws.Range("F25").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("D26").Value = "'0.9241"
$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'0.9239"
$ws.Range("F26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0.9742"
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'0.9749"
$ws.Range("F27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# --- Column A width adjustment ---
$ws.Columns.Item(1).ColumnWidth = 32.61

# --- Update active selection to D12 (matches author's final cursor position) ---
$ws.Range("D12").Select()

$excel.CutCopyMode = $false
